$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.828.00"
$ws.Range("E2").Value = "  +3.17%  "
$ws.Range("D3").Value = "3.132.02"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "3.124.51"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").Value = "  +19.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.77%  "
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "3.650.17"
$ws.Range("E16").Value = "  +2.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "63.757.01"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").Value = "3.129.68"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.40%  "
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").Value = "0.0₃0872"
$ws.Range("E34").Value = "  +7.78%  "
$ws.Range("E35").Value = "  +10.07%  "
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "452.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0373"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "2.920.23"
$ws.Range("E43").Value = "  +5.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.278"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.35%  "
$ws.Range("E45").Value = "  +3.25%  "
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.41%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.69"
$ws.Range("D51").Style = "Normal"
